# The author selected row 8's "zdjecie"/"podpis" cells (G8:H8) and deleted
# their contents, then scrolled the view down so row 7 becomes the new
# top-left visible row, leaving G8 as the active/selected cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the photo filename / caption text that used to live in G8 and H8.
$ws.Range("G8:H8").ClearContents()

# Scroll the window so row 7 is at the top, matching the saved view state.
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1

# Leave the selection on G8, as recorded in the saved file.
$ws.Range("G8").Select()
